$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Agency" column (H) ------------------------------------------
# Maps each project row to its funding-agency category.
$agency = @{
    1  = "Agency"
    2  = "DST"
    3  = "DBT"
    4  = "SERB"
    5  = "ICMR"
    6  = "DBT"
    7  = "SERB"
    8  = "DBT"
    9  = "DBT"
    10 = "SERB"
    11 = "SERB"
    12 = "DST"
    13 = "DST"
    14 = "Industry"
    15 = "IIT Hyderabad"
    16 = "DBT"
    17 = "ICMR"
    18 = "IIT Hyderabad"
    19 = "DBT"
    20 = "DBT"
    21 = "SERB"
    22 = "MOE"
    23 = "ICMR"
    24 = "ICMR"
    25 = "SERB"
    26 = "International"
    27 = "ICMR"
    28 = "SERB"
    29 = "International"
    30 = "SERB"
    31 = "Other sources"
    32 = "DBT"
    33 = "DBT"
    34 = "DBT"
    35 = "DBT"
}

# Write the rows that introduce a brand-new label first (and in the exact
# order those labels were first typed), so new shared-string entries land in
# the same sequence as the authored workbook; then fill in the rest.
$firstRows = @(1, 2, 14, 15, 31, 22, 26)
foreach ($row in $firstRows) {
    $ws.Cells.Item($row, 8).Value = $agency[$row]
}
for ($row = 1; $row -le 35; $row++) {
    if ($firstRows -notcontains $row) {
        $ws.Cells.Item($row, 8).Value = $agency[$row]
    }
}

# --- Column widths (B -> 105 chars, D -> 65 chars, best-fit) ---------------
$ws.Columns.Item(2).ColumnWidth = 104.16666666666667
$ws.Columns.Item(4).ColumnWidth = 64.16666666666667

# --- AutoFilter on the Funding Agency column (D) ---------------------------
$ws.Range("D1:D35").AutoFilter() | Out-Null

$filterName = $ws.Names.Add("_FilterDatabase", "=Sheet1!`$D`$1:`$D`$35")
$filterName.Visible = $false

# --- Final selection on H35 -------------------------------------------------
$ws.Range("H35").Select() | Out-Null
